# Cap nhat 10 bai toan thu nghiem
# Update experiment results for the "problem_stats" sheet: refresh the
# Problem 5, Problem 8 and Problem 9 rows with newly measured values, then
# refresh the derived Min / Max / Mean / Std / Std-percent summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Problem 5 (row 6) ---
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 91.67
$ws.Range("D6").Value = 99.58
$ws.Range("E6").Value = 136
$ws.Range("F6").Value = 40.15
$ws.Range("G6").Value = 86
$ws.Range("H6").Value = 39.39

# --- Problem 8 (row 9) ---
$ws.Range("B9").Value = 15
$ws.Range("C9").Value = 6.67
$ws.Range("D9").Value = 51.55
$ws.Range("E9").Value = 275
$ws.Range("F9").Value = 85.16

# --- Problem 9 (row 10) ---
$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 13.33
$ws.Range("D10").Value = 58.46
$ws.Range("E10").Value = 205.5
$ws.Range("F10").Value = 144.01

# --- Min (row 12) ---
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = 6.67
$ws.Range("D12").Value = 51.55

# --- Max (row 13) ---
$ws.Range("E13").Value = 275

# --- Mean (row 14) ---
$ws.Range("B14").Value = 11.4
$ws.Range("C14").Value = 69.64500000000001
$ws.Range("D14").Value = 80.392
$ws.Range("E14").Value = 146.124
$ws.Range("F14").Value = 79.672
$ws.Range("G14").Value = 91.45999999999999
$ws.Range("H14").Value = 47.806

# --- Std (row 15) ---
$ws.Range("B15").Value = 5.015531433014408
$ws.Range("C15").Value = 33.37109902295698
$ws.Range("D15").Value = 16.56148128103951
$ws.Range("E15").Value = 60.43630148989742
$ws.Range("F15").Value = 49.46316312839957
$ws.Range("G15").Value = 9.991706560943431
$ws.Range("H15").Value = 24.43989325308566

# --- Std / (max - min) % (row 16) ---
$ws.Range("B16").Value = 33.43687622009605
$ws.Range("C16").Value = 35.75602595409512
$ws.Range("D16").Value = 34.18262390307433
$ws.Range("E16").Value = 29.48112267799874
$ws.Range("F16").Value = 30.71101647112851
$ws.Range("G16").Value = 33.3056885364781
$ws.Range("H16").Value = 33.75675863685865
